$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right 5 -> 4, Wrong -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right 95 -> 76, Wrong -2 -> -4, summary text "95 / 140" -> "72 / 112"
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "72 / 112"
